$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# ---- Header row (row 1) ----
$headers = @(
    "Date",
    "Model Name",
    "Exact Precision (Micro Avg)",
    "Exact Recall (Micro Avg)",
    "Exact F1 Score (Micro Avg)",
    "Exact Precision (Macro Avg)",
    "Exact Recall (Macro Avg)",
    "Exact F1 Score (Macro Avg)",
    "Exact Precision (Weighted Avg)",
    "Exact Recall (Weighted Avg)",
    "Exact F1 Score (Weighted Avg)",
    "Partial Precision",
    "Partial Recall",
    "Partial F1 Score",
    "Partial TP",
    "Partial FP",
    "Partial FN",
    "Support",
    "Accuracy",
    "Result Link",
    "Stats Link",
    "No of GPU Used",
    "Power Consumption"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Bold, thin-bordered, center/top aligned header style
$headerRange = $ws.Range("A1:W1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = "Continuous"
$headerRange.Borders.Weight = "Thin"
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# ---- Data row (row 2) ----
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "09/11/2025"
$ws.Cells.Item(2, 1).Style = "Normal"

$ws.Cells.Item(2, 2).Value = "Llama-3.1-8B-Instruct"
$ws.Cells.Item(2, 3).Value = 0.3902439024390244
$ws.Cells.Item(2, 4).Value = 0.2730375426621161
$ws.Cells.Item(2, 5).Value = 0.321285140562249
$ws.Cells.Item(2, 6).Value = 0.2482302011713776
$ws.Cells.Item(2, 7).Value = 0.1229980723685647
$ws.Cells.Item(2, 8).Value = 0.1494622268994377
$ws.Cells.Item(2, 9).Value = 0.5222522642317864
$ws.Cells.Item(2, 10).Value = 0.2730375426621161
$ws.Cells.Item(2, 11).Value = 0.3377222915882462
$ws.Cells.Item(2, 12).Value = 0.4876847290640394
$ws.Cells.Item(2, 13).Value = 0.339041095890411
$ws.Cells.Item(2, 14).Value = 0.4
$ws.Cells.Item(2, 15).Value = 99
$ws.Cells.Item(2, 16).Value = 104
$ws.Cells.Item(2, 17).Value = 193
$ws.Cells.Item(2, 18).Value = 293
$ws.Cells.Item(2, 19).Value = 0.9476213951195465
$ws.Cells.Item(2, 20).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Llama-3.1-8B-Instruct_3_shot.txt"
$ws.Cells.Item(2, 21).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Llama-3.1-8B-Instruct_3_shot.txt"
$ws.Cells.Item(2, 22).Value = "4 MLGPU"
$ws.Cells.Item(2, 23).Value = "0.014 kWh"
$ws.Cells.Item(2, 24).Value = 562
